# Re-implements the pptx table-writer change described in the commit:
# "Using internal column widths in pptx writer tables" — the two
# <a:gridCol> widths in the slide's table move from 2501900 EMU
# (197pt) to 2514600 EMU (198pt) each, and the paragraphs whose
# <a:pPr> get rewritten as part of that pass end up with their
# lvl/indent/marL attributes re-ordered.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shapes we need: the title placeholder, the table's
# graphic frame, and the free-floating caption textbox.
$titleShape = $null
$tableShape = $null
$textboxShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
    } elseif ($shp.Name -eq "Title 1") {
        $titleShape = $shp
    } elseif ($shp.Name -eq "TextBox 3") {
        $textboxShape = $shp
    }
}

# --- Table: apply the new internal column widths (198pt = 2514600 EMU) ---
$tbl = $tableShape.Table
$tbl.Columns.Item(1).Width = 198
$tbl.Columns.Item(2).Width = 198

# Touch each table-cell paragraph's formatting so it is re-written
# with the writer's current attribute ordering (lvl/indent/marL).
$tbl.Cell(1, 1).Shape.TextFrame.TextRange.IndentLevel = 1
$tbl.Cell(1, 2).Shape.TextFrame.TextRange.IndentLevel = 1

# --- Title placeholder paragraph formatting re-write ---
$titleShape.TextFrame.TextRange.IndentLevel = 1

# --- Caption textbox paragraph formatting re-write ---
$textboxShape.TextFrame.TextRange.IndentLevel = 1
